$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 2).Value = 180
$ws.Cells.Item(2, 4).Value = 300
$ws.Cells.Item(2, 6).Value = 45
$ws.Cells.Item(3, 3).Value = 0.0
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(5, 3).Value = 0.5555555555555556
$ws.Cells.Item(6, 3).Value = 0.5555555555555556
$ws.Cells.Item(7, 3).Value = 0.5555555555555556
$ws.Cells.Item(8, 3).Value = 0.5555555555555556
$ws.Cells.Item(9, 3).Value = 0.5555555555555556
$ws.Cells.Item(10, 3).Value = 0.5555555555555556
$ws.Cells.Item(11, 3).Value = 1.111111111111111
$ws.Cells.Item(12, 3).Value = 1.111111111111111
$ws.Cells.Item(13, 3).Value = 1.111111111111111
$ws.Cells.Item(14, 3).Value = 1.111111111111111
$ws.Cells.Item(15, 3).Value = 1.111111111111111
$ws.Cells.Item(16, 3).Value = 1.111111111111111
$ws.Cells.Item(17, 3).Value = 1.666666666666667
$ws.Cells.Item(18, 3).Value = 1.666666666666667
$ws.Cells.Item(19, 3).Value = 1.666666666666667
$ws.Cells.Item(20, 3).Value = 1.666666666666667
$ws.Cells.Item(21, 3).Value = 1.666666666666667
$ws.Cells.Item(22, 3).Value = 1.666666666666667
$ws.Cells.Item(23, 3).Value = 2.222222222222222
$ws.Cells.Item(24, 3).Value = 2.222222222222222
$ws.Cells.Item(25, 3).Value = 2.222222222222222
$ws.Cells.Item(26, 3).Value = 2.222222222222222
$ws.Cells.Item(27, 3).Value = 2.222222222222222
$ws.Cells.Item(28, 3).Value = 2.222222222222222
$ws.Cells.Item(29, 3).Value = 2.777777777777778
$ws.Cells.Item(30, 3).Value = 2.777777777777778
$ws.Cells.Item(31, 3).Value = 2.777777777777778
$ws.Cells.Item(32, 3).Value = 2.777777777777778
$ws.Cells.Item(33, 3).Value = 2.777777777777778
$ws.Cells.Item(34, 3).Value = 2.777777777777778
$ws.Cells.Item(35, 3).Value = 3.333333333333333
$ws.Cells.Item(36, 3).Value = 3.333333333333333
$ws.Cells.Item(37, 3).Value = 3.333333333333333
$ws.Cells.Item(38, 3).Value = 3.333333333333333
$ws.Cells.Item(39, 3).Value = 3.333333333333333
